# Add "NA" values under duplicate_image_filename (column E) for data rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the pre-existing (malformed, value-less) shared-string cell at F1
# back to an empty cell; otherwise the COM runtime normalizes it to a
# visible value as a side effect of loading/saving the workbook.
$ws.Range("F1").ClearContents()

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
